$d = $word.ActiveDocument

$replacements = @(
    @{old="61×33="; new="32×14="},
    @{old="45×59="; new="24×90="},
    @{old="66×61="; new="19×41="},
    @{old="83×23="; new="42×47="},
    @{old="75×32="; new="80×38="},
    @{old="22×84="; new="66×70="},
    @{old="51×13="; new="76×18="},
    @{old="23×90="; new="25×75="},
    @{old="87×21="; new="66×95="},
    @{old="52×57="; new="82×61="},
    @{old="82×11="; new="29×87="},
    @{old="25×59="; new="57×94="},
    @{old="47×99="; new="75×75="},
    @{old="48×16="; new="20×28="},
    @{old="20×59="; new="65×81="},
    @{old="75×87="; new="18×13="},
    @{old="67×22="; new="27×62="},
    @{old="14×70="; new="97×49="},
    @{old="74×57="; new="86×67="},
    @{old="27×23="; new="96×29="},
    @{old="41×93="; new="33×11="},
    @{old="96×41="; new="55×66="},
    @{old="26×89="; new="87×28="},
    @{old="50×37="; new="97×19="},
    @{old="85×54="; new="49×85="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
